$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: Kelly Oubre Jr. -> Stephen Curry
$ws.Range("A5").Value = "Stephen Curry"
$ws.Range("B5").Value = "PG,SG"
$ws.Range("C5").Value = "Golden State Warriors"

# Row 14: Daniel Gafford -> Kelly Oubre Jr.
$ws.Range("A14").Value = "Kelly Oubre Jr."
$ws.Range("B14").Value = "SG,SF"
$ws.Range("C14").Value = "Philadelphia 76ers"

# Row 15: Keegan Murray -> Daniel Gafford
$ws.Range("A15").Value = "Daniel Gafford"
$ws.Range("B15").Value = "PF,C"
$ws.Range("C15").Value = "Dallas Mavericks"

# Row 16: Stephen Curry -> Keegan Murray
$ws.Range("A16").Value = "Keegan Murray"
$ws.Range("B16").Value = "SF,PF"
$ws.Range("C16").Value = "Sacramento Kings"
